$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new cell values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the bold / bordered / centered-top style on B1 ...
$r1 = $ws.Range("B1")
$r1.HorizontalAlignment = -4108   # xlCenter
$r1.VerticalAlignment = -4160     # xlTop
$r1.Font.Bold = $true
$r1.Borders.LineStyle = 1         # xlContinuous (thin)

# ... then replicate the exact same style onto A2 via a format copy,
# so both cells end up referencing the very same cell style (xf).
$r1.Copy()
$ws.Range("A2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
